# Applies the updated cryptos list values to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "69.711.96", "1.00")
# that must stay as literal text, matching the source data exactly -
# including trailing zeros and the dotted "thousands" groupings used
# by this feed. Temporarily force the column to text format so Excel
# does not reinterpret / round the strings as numbers, then restore
# the column back to its normal (General) format.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "69.711.96"
$ws.Range("E2").Value = "  -0.79%  "

$ws.Range("D3").Value = "3.555.75"
$ws.Range("E3").Value = "  -1.46%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "197.16"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").Value = "587.30"
$ws.Range("E6").Value = "  -2.83%  "

$ws.Range("E7").Value = "  -2.22%  "

$ws.Range("D8").Value = "1.00"

$ws.Range("E9").Value = "  +0.58%  "

$ws.Range("D10").Value = "0.633"
$ws.Range("E10").Value = "  -2.41%  "

$ws.Range("D11").Value = "52.89"
$ws.Range("E11").Value = "  -1.86%  "

$ws.Range("D12").Value = "0.0000289"
$ws.Range("E12").Value = "  -4.84%  "

$ws.Range("D13").Value = "9.28"
$ws.Range("E13").Value = "  -3.11%  "

$ws.Range("D14").Value = "4.116.82"
$ws.Range("E14").Value = "  -1.49%  "

$ws.Range("D15").Value = "671.06"
$ws.Range("E15").Value = "  +12.12%  "

$ws.Range("D16").Value = "69.788.42"
$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("D17").Value = "3.563.24"
$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D18").Value = "12.54"
$ws.Range("E18").Value = "  -4.38%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "18.51"
$ws.Range("E19").Value = "  -3.10%  "

$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "0.121"
$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("D21").Value = "0.968"
$ws.Range("E21").Value = "  -2.88%  "

$ws.Range("D22").Value = "18.08"
$ws.Range("E22").Value = "  +1.49%  "

$ws.Range("D23").Value = "5.37"
$ws.Range("E23").Value = "  +3.60%  "

$ws.Range("D24").Value = "105.64"
$ws.Range("E24").Value = "  +3.45%  "

$ws.Range("D25").Value = "4.40"
$ws.Range("E25").Value = "  -4.82%  "

$ws.Range("E26").Value = "  -3.27%  "

$ws.Range("D27").Value = "10.20"
$ws.Range("E27").Value = "  -5.25%  "

$ws.Range("D28").Value = "9.65"
$ws.Range("E28").Value = "  +0.28%  "

$ws.Range("D29").Value = "33.56"
$ws.Range("E29").Value = "  -0.88%  "

$ws.Range("D30").Value = "4.40"
$ws.Range("E30").Value = "  -7.80%  "

$ws.Range("D31").Value = "6.79"
$ws.Range("E31").Value = "  -5.30%  "

$ws.Range("D32").Value = "11.78"
$ws.Range("E32").Value = "  -4.15%  "

$ws.Range("D33").Value = "0.112"
$ws.Range("E33").Value = "  -4.43%  "

$ws.Range("D34").Value = "62.12"
$ws.Range("E34").Value = "  -1.85%  "

$ws.Range("D35").Value = "3.781.59"
$ws.Range("E35").Value = "  -3.29%  "

$ws.Range("E36").Value = "  -8.59%  "

$ws.Range("E37").Value = "  +5.98%  "

$ws.Range("E38").Value = "  -0.06%  "

$ws.Range("D39").Value = "502.25"
$ws.Range("E39").Value = "  -3.48%  "

$ws.Range("E40").Value = "  -6.32%  "

$ws.Range("D41").Value = "0.372"
$ws.Range("E41").Value = "  -4.48%  "

$ws.Range("E42").Value = "  +1.18%  "

$ws.Range("D43").Value = "34.69"
$ws.Range("E43").Value = "  -6.13%  "

$ws.Range("E44").Value = "  -0.83%  "

$ws.Range("D45").Value = "2.87"
$ws.Range("E45").Value = "  +0.63%  "

$ws.Range("E46").Value = "  -1.50%  "

$ws.Range("E47").Value = "  -2.40%  "

$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("D49").Value = "8.37"
$ws.Range("E49").Value = "  -3.04%  "

$ws.Range("D50").Value = "1.77"
$ws.Range("E50").Value = "  +19.32%  "

$ws.Range("D51").Value = "2.72"
$ws.Range("E51").Value = "  +62.60%  "

$priceCol.NumberFormat = "general"
